# [TMS][Update] Add line manager to db and import user function
#
# Adds a new "Line Manager" column (AH) to the Staff List sheet:
#   - AH1:AH3 merged header cell "Line Manager" styled like the adjoining
#     AG1:AG3 ("Skype") header (same fill/border/font treatment).
#   - AH4 continues the column-index numbering row (32 -> 33).
#   - AH5 is left blank (the manager row itself has no line manager).
#   - AH6 / AH7 contain the line manager's e-mail address as a mailto
#     hyperlink, matching the styling already used for the AF (Email)
#     column's hyperlink cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header (row 1-3), merged AH1:AH3, formatted like AG1:AG3 -------------
$ws.Range("AG1").Copy()
$ws.Range("AH1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AG2").Copy()
$ws.Range("AH2").PasteSpecial(-4122)

$ws.Range("AG3").Copy()
$ws.Range("AH3").PasteSpecial(-4122)

$ws.Range("AH1").Value = "Line Manager"
$ws.Range("AH1:AH3").Merge()

# --- Column width (approx. 32.43 chars wide, bestFit like its neighbours) -
$ws.Columns("AH").ColumnWidth = 31.67

# --- Row 4: sequence number continues --------------------------------------
$ws.Range("AG4").Copy()
$ws.Range("AH4").PasteSpecial(-4122)
$ws.Range("AH4").Value = 33

# --- Row 5: the manager's own row - no line manager, cell left blank -------
$ws.Range("AF6").Copy()
$ws.Range("AH5").PasteSpecial(-4122)

# --- Rows 6-7: line manager e-mail + hyperlink ------------------------------
$ws.Range("AF5").Copy()
$ws.Range("AH6").PasteSpecial(-4122)
$ws.Range("AH6").Value = "quanghuy@easia-travel.com"
$ws.Hyperlinks.Add($ws.Range("AH6"), "mailto:quanghuy@easia-travel.com")

$ws.Range("AF5").Copy()
$ws.Range("AH7").PasteSpecial(-4122)
$ws.Range("AH7").Value = "quanghuy@easia-travel.com"
$ws.Hyperlinks.Add($ws.Range("AH7"), "mailto:quanghuy@easia-travel.com")

# --- View state: scroll right and land selection on AH6 --------------------
$win = $excel.ActiveWindow
$win.ScrollColumn = 31
$win.ScrollRow = 1
$ws.Range("AH6").Select()

$excel.CutCopyMode = 0
